$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 2532
$ws.Range("B2").Value = "Pietra Sales"
$ws.Range("C2").Value = "TI"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45100
$ws.Range("G2").Value = 4180.35

# Row 3
$ws.Range("A3").Value = 71473
$ws.Range("B3").Value = "Rael Pereira"
$ws.Range("C3").Value = "Recursos Humanos"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45079
$ws.Range("G3").Value = 4105.87

# Row 4
$ws.Range("A4").Value = 26437
$ws.Range("B4").Value = "Stephany Silva"
$ws.Range("C4").Value = "Recursos Humanos"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 45085
$ws.Range("G4").Value = 2073.68

# Row 5
$ws.Range("A5").Value = 9289
$ws.Range("B5").Value = "Valentina Costela"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("D5").Value = "Consulta medica"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 45087
$ws.Range("G5").Value = 3371.66

# Row 6
$ws.Range("A6").Value = 18061
$ws.Range("B6").Value = "Helena Novais"
$ws.Range("C6").Value = "Atendimento ao Cliente"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45101
$ws.Range("G6").Value = 4300.22

# Row 7
$ws.Range("A7").Value = 91687
$ws.Range("B7").Value = "Eloá Silva"
$ws.Range("C7").Value = "Recursos Humanos"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45089
$ws.Range("G7").Value = 6890.23

# Row 8
$ws.Range("A8").Value = 2838
$ws.Range("B8").Value = "Sr. Benicio Silva"
$ws.Range("C8").Value = "Operacoes"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45104
$ws.Range("G8").Value = 4238.11

# Row 9
$ws.Range("A9").Value = 82708
$ws.Range("B9").Value = "Dr. Cauã Albuquerque"
$ws.Range("C9").Value = "Atendimento ao Cliente"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45101
$ws.Range("G9").Value = 8259.48

# Row 10
$ws.Range("A10").Value = 78451
$ws.Range("B10").Value = "Dom da Mota"
$ws.Range("C10").Value = "TI"
$ws.Range("D10").Value = "Doenca"
$ws.Range("E10").Value = 1
$ws.Range("G10").Value = 5899.55

# Row 11
$ws.Range("A11").Value = 78784
$ws.Range("B11").Value = "Isis da Paz"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 8
$ws.Range("G11").Value = 2888.7
